$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = ' depreciation and amortization'
$ws.Range("B2").ClearContents()
$ws.Range("C2").ClearContents()
$ws.Range("D2").ClearContents()

$ws.Range("A3").Value = ' stock-based compensation expense'
$ws.Range("B3").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("D3").ClearContents()

$ws.Range("A4").Value = ' tax benefit from stock-based compensation'
$ws.Range("B4").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("D4").ClearContents()

$ws.Range("A5").Value = ' other'
$ws.Range("B5").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("D5").ClearContents()

$ws.Range("A6").Value = ' trade and other accounts receivable net'
$ws.Range("B6").ClearContents()
$ws.Range("C6").ClearContents()
$ws.Range("D6").ClearContents()

$ws.Range("A7").Value = ' inventories'
$ws.Range("B7").ClearContents()
$ws.Range("C7").ClearContents()
$ws.Range("D7").ClearContents()

$ws.Range("A8").Value = ' prepaid expenses and other'
$ws.Range("B8").ClearContents()
$ws.Range("C8").ClearContents()
$ws.Range("D8").ClearContents()

$ws.Range("A9").Value = ' deposits and other'
$ws.Range("B9").ClearContents()
$ws.Range("C9").ClearContents()
$ws.Range("D9").ClearContents()

$ws.Range("A10").Value = ' accounts payable'
$ws.Range("B10").ClearContents()
$ws.Range("C10").ClearContents()
$ws.Range("D10").ClearContents()

$ws.Range("A11").Value = ' accrued expenses'
$ws.Range("B11").ClearContents()
$ws.Range("C11").ClearContents()
$ws.Range("D11").ClearContents()

$ws.Range("A12").Value = ' deferred rent'
$ws.Range("B12").ClearContents()
$ws.Range("C12").ClearContents()
$ws.Range("D12").ClearContents()

$ws.Range("A13").Value = ' other net long-term cash provided liabilities by operating activities'
$ws.Range("B13").ClearContents()
$ws.Range("C13").ClearContents()
$ws.Range("D13").ClearContents()

$ws.Range("A14").Value = ' additions to property and equipment'
$ws.Range("B14").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Range("D14").ClearContents()

$ws.Range("A15").Value = ' proceeds net from cash sale-leaseback used in investing transactions activities'
$ws.Range("B15").ClearContents()
$ws.Range("C15").ClearContents()
$ws.Range("D15").ClearContents()

$ws.Range("A16").Value = ' exercise of employee stock options'
$ws.Range("B16").ClearContents()
$ws.Range("C16").ClearContents()
$ws.Range("D16").ClearContents()

$ws.Range("A17").Value = ' tax benefit from stock-based compensation'
$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()
$ws.Range("D17").ClearContents()

$ws.Range("A18").Value = ' proceeds net from cash issuance used in financing of common activities stock under employee benefit plans'
$ws.Range("B18").ClearContents()
$ws.Range("C18").ClearContents()
$ws.Range("D18").ClearContents()

$ws.Range("A19").Value = ' net increase (decrease) in cash and cash equivalents'
$ws.Range("B19").ClearContents()
$ws.Range("C19").ClearContents()
$ws.Range("D19").ClearContents()

$ws.Range("A20").Value = ' cash and cash equivalents at beginning of period'
$ws.Range("B20").ClearContents()
$ws.Range("C20").ClearContents()
$ws.Range("D20").ClearContents()

$ws.Range("A21").Value = ' cash and cash equivalents at end of the period accompanying notes are an integral part of the consolidated financial statements.'
$ws.Range("B21").ClearContents()
$ws.Range("C21").ClearContents()
$ws.Range("D21").ClearContents()

$ws.Range("A22").Value = ' business combinations net of cash acquired'
$ws.Range("B22").Value = -13
$ws.Range("C22").Value = -45
$ws.Range("D22").Value = -18

$ws.Range("A23").Value = ' net cash used in investing activities'
$ws.Range("B23").Value = -3132
$ws.Range("C23").Value = -1436
$ws.Range("D23").Value = -2337

$ws.Range("A24").Value = ' proceeds from issuances of convertible and other debt'
$ws.Range("B24").Value = 9713
$ws.Range("C24").Value = 10669
$ws.Range("D24").Value = 6176

$ws.Range("A25").Value = ' repayments of convertible and other debt'
$ws.Range("B25").Value = -11623
$ws.Range("C25").Value = -9161
$ws.Range("D25").Value = -5247

$ws.Range("A26").Value = ' collateralized lease repayments'
$ws.Range("B26").Value = -240
$ws.Range("C26").Value = -389
$ws.Range("D26").Value = -559

$ws.Range("A27").Value = ' proceeds from exercises of stock options and other stock issuances'
$ws.Range("B27").Value = 417
$ws.Range("C27").Value = 263
$ws.Range("D27").Value = 296

$ws.Range("A28").Value = ' principal payments on finance leases'
$ws.Range("B28").Value = -338
$ws.Range("C28").Value = -321
$ws.Range("D28").Value = -181

$ws.Range("A29").Value = ' debt issuance costs'
$ws.Range("B29").Value = -6
$ws.Range("C29").Value = -37
$ws.Range("D29").Value = -15

$ws.Range("A30").Value = ' proceeds from investments by noncontrolling interests in subsidiaries'
$ws.Range("B30").Value = 24
$ws.Range("C30").Value = 279
$ws.Range("D30").Value = 437

$ws.Range("A31").Value = ' distributions paid to noncontrolling interests in subsidiaries gl)'
$ws.Range("B31").Value = -208
$ws.Range("C31").Value = 0
$ws.Range("D31").Value = -227

$ws.Range("A32").Value = ' payments for buy-outs of noncontrolling interests in subsidiaries'
$ws.Range("B32").Value = -35
$ws.Range("C32").Value = -9
$ws.Range("D32").Value = -6

$ws.Range("A33").Value = ' net cash provided by financing activities'
$ws.Range("B33").Value = 9973
$ws.Range("C33").Value = 1529
$ws.Range("D33").Value = 574

$ws.Range("A34").Value = ' effect of exchange rate changes on cash and cash equivalents and restricted cash'
$ws.Range("B34").Value = 334
$ws.Range("C34").Value = 8
$ws.Range("D34").Value = -23

$ws.Range("A35").Value = ' net increase in cash and cash equivalents and restricted cash'
$ws.Range("B35").Value = 13118
$ws.Range("C35").Value = 2506
$ws.Range("D35").Value = 312

$ws.Range("A36").Value = ' cash and cash equivalents and restricted cash beginning of period'
$ws.Range("B36").Value = 6783
$ws.Range("C36").Value = 4277
$ws.Range("D36").Value = 3965

$ws.Range("A37").Value = ' cash and cash equivalents and restricted cash end of period'
$ws.Range("B37").Value = 19901
$ws.Range("C37").Value = 6783
$ws.Range("D37").Value = 4277

$ws.Range("A38").Value = ' acquisitions of property and equipment included in liabilities'
$ws.Range("B38").Value = 1088
$ws.Range("C38").Value = 562
$ws.Range("D38").Value = 249
